# correção nos dados e inicio da analise PNAD 2009
#
# The original sheet had two "category header" rows (row 5 "situação do
# domicílio" and row 8 "grandes regiões e unidades da federação") that only
# carried a label in column A with no data beside them. This edit removes
# those two empty/placeholder rows (shifting everything below them up), and
# fixes the row-2 header labels so that the "total" columns (B and F) read
# "total" instead of the stray "unnamed: 1_level_1" / "unnamed: 5_level_1"
# pandas-export artifacts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so row 5's index isn't affected.
$ws.Rows("8").Delete()
$ws.Rows("5").Delete()

# Correct the header row: B2/F2 should match C2's "total" label.
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
